$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new weekly price rows for Cereza (cherry) that need to be
# inserted above the existing "row 167" data, pushing the old rows 167-174 down
# to become rows 169-176 (one of which also gets a date correction).

# Insert two blank rows at position 167 (run twice so both land above the old data).
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# --- New row 167: Early Burlat / Primera, Región de O'Higgins ---
$ws.Range("A167").Value = 9
$ws.Range("B167").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C167").Value = "Metropolitana"
$ws.Range("D167").Value = 44516
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100103
$ws.Range("H167").Value = "Frutos de hueso (carozo)"
$ws.Range("I167").Value = 100103001
$ws.Range("J167").Value = "Cereza"
$ws.Range("K167").Value = "Early Burlat"
$ws.Range("L167").Value = "Primera"
$ws.Range("M167").Value = 450
$ws.Range("N167").Value = 30000
$ws.Range("O167").Value = 30000
$ws.Range("P167").Value = 30000
$ws.Range("Q167").Value = "$/bandeja 10 kilos"
$ws.Range("R167").Value = "Región de O'Higgins"
$ws.Range("S167").Value = 3000
$ws.Range("T167").Value = 10

# --- New row 168: Royal Dawn / Primera, Provincia de Curicó ---
$ws.Range("A168").Value = 9
$ws.Range("B168").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C168").Value = "Metropolitana"
$ws.Range("D168").Value = 44516
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100103
$ws.Range("H168").Value = "Frutos de hueso (carozo)"
$ws.Range("I168").Value = 100103001
$ws.Range("J168").Value = "Cereza"
$ws.Range("K168").Value = "Royal Dawn"
$ws.Range("L168").Value = "Primera"
$ws.Range("M168").Value = 300
$ws.Range("N168").Value = 24000
$ws.Range("O168").Value = 24000
$ws.Range("P168").Value = 24000
$ws.Range("Q168").Value = "$/bandeja 8 kilos"
$ws.Range("R168").Value = "Provincia de Curicó"
$ws.Range("S168").Value = 3000
$ws.Range("T168").Value = 8

# The old row 172 (Santina / Primera / Región de O'Higgins), now shifted to row
# 174, had its date corrected from 44179 to 44217.
$ws.Range("D174").Value = 44217
